# Update the cryptos list (Price / Volume(1h) columns) with the latest
# scraped snapshot. Values that look numeric (single decimal point) are
# prefixed with an apostrophe so Excel stores them as text, matching the
# original inline-string cell type instead of converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.529.68"
$ws.Range("E2").Value = "  +4.78%  "
$ws.Range("D3").Value = "1.591.68"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "'214.39"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "'23.80"
$ws.Range("E8").Value = "  +7.83%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'0.0887"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "1.819.92"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "1.588.96"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'0.530"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "28.513.72"
$ws.Range("E16").Value = "  +4.90%  "
$ws.Range("D17").Value = "'63.99"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "'233.10"
$ws.Range("E18").Value = "  +7.81%  "
$ws.Range("D19").Value = "'7.53"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "0.0₃0709"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'151.87"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "'15.35"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "1.422.48"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +8.95%  "
$ws.Range("D40").Value = "'0.546"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("D45").Value = "'0.974"
$ws.Range("E45").Value = "  -3.02%  "
$ws.Range("D46").Value = "'64.65"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "1.729.55"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "'87.69"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'39.74"
$ws.Range("E51").Value = "  +17.18%  "
